# Add new columns I (I0) and J (IF) to the sheet, matching the
# existing header style used by the other header cells (e.g. H1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new header cells I1 and J1.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting of the existing header cell (H1) onto the new
# header cells so they share the same bold/centered/bordered style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-15 for columns I and J.
$data = @{
    2  = @(6, 8)
    3  = @(4, 5)
    4  = @(1, 2)
    5  = @(8, 8)
    6  = @(4, 5)
    7  = @(1, 3)
    8  = @(1, 5)
    9  = @(1, 5)
    10 = @(1, 4)
    11 = @(1, 4)
    12 = @(1, 5)
    13 = @(1, 4)
    14 = @(1, 3)
    15 = @(1, 2)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 9).Value = $values[0]
    $ws.Cells.Item($row, 10).Value = $values[1]
}
